$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D; existing D:K shift to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting from column F (old column D) onto new columns D:E
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the two new columns (D, E) with the newest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 27600
$ws.Range("E8").Value = 32700
$ws.Range("D9").Value = 9400
$ws.Range("E9").Value = 8800
$ws.Range("D10").Value = 18200
$ws.Range("E10").Value = 23900
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 14200
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 10400
$ws.Range("E15").Value = 10900
$ws.Range("D17").Value = 37500
$ws.Range("E17").Value = 23100
$ws.Range("D18").Value = -9900
$ws.Range("E18").Value = 9600
$ws.Range("D20").Value = 1800
$ws.Range("E20").Value = -3300
$ws.Range("D21").Value = 2300
$ws.Range("E21").Value = 17300
$ws.Range("D22").Value = 300
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -8500
$ws.Range("E23").Value = 6300
$ws.Range("D24").Value = -1400
$ws.Range("E24").Value = 600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -7100
$ws.Range("E26").Value = 5700
$ws.Range("D27").Value = -7100
$ws.Range("E27").Value = 5700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1800
$ws.Range("E32").Value = 3300
$ws.Range("D33").Value = -7100
$ws.Range("E33").Value = 5700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -7100
$ws.Range("E35").Value = 5700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 3400
$ws.Range("E41").Value = 3800
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 13200
$ws.Range("E43").Value = 13500
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 300
$ws.Range("E45").Value = 400
$ws.Range("D46").Value = 16800
$ws.Range("E46").Value = 17600
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 542000
$ws.Range("E48").Value = 509100
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 8200
$ws.Range("E52").Value = 6900
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 567100
$ws.Range("E54").Value = 533700
$ws.Range("D57").Value = 51900
$ws.Range("E57").Value = 42500
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("E59").Value = 6400
$ws.Range("D60").Value = 51900
$ws.Range("E60").Value = 48900
$ws.Range("D61").Value = 39500
$ws.Range("E61").Value = 17000
$ws.Range("D62").Value = 13100
$ws.Range("E62").Value = 10200
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 104500
$ws.Range("E66").Value = 76100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -32400
$ws.Range("E72").Value = -25300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 462600
$ws.Range("E76").Value = 457600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -7100
$ws.Range("E81").Value = 5700
$ws.Range("D83").Value = 10400
$ws.Range("E83").Value = 10900
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 18400
$ws.Range("E89").Value = 18700
$ws.Range("D91").Value = -600
$ws.Range("E91").Value = -800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -41400
$ws.Range("E94").Value = -45400
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 22600
$ws.Range("E100").Value = 17000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -400
$ws.Range("E102").Value = -9700

# Correct the handful of cells in the shifted F:M range whose values changed
# (not a pure positional shift of the old D:K data)
$ws.Range("H24").Value = 0
$ws.Range("H26").Value = 2400
$ws.Range("H27").Value = 2400
$ws.Range("H29").Value = -7000
$ws.Range("F91").Value = -2200
$ws.Range("G91").Value = -1100
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = -100
$ws.Range("J91").Value = 2600
